# Fix auto-increase enrollment flags and add debug logging for participation tracking
#
# Updates the projection results for rows 2-6 (plan years 1-5) on the active
# worksheet to reflect corrected "Participating" headcounts (auto-increase
# enrollment now flows into participation correctly) and the downstream
# metrics derived from it (rates, contributions, totals, costs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column letter -> new value
$updates = @{
    2 = @{
        'C' = 9463
        'D' = 8393
        'E' = 0.8869280355067104
        'F' = 0.8855243722304283
        'G' = 0.09671184592974379
        'H' = 0.08564069665418225
        'I' = 41239758.28268903
        'J' = 14425047.08800052
        'L' = 14425047.08800052
        'M' = 55664805.37068957
        'N' = 800750889.2872001
        'O' = 783051082.2832
        'P' = 0.01801440033471731
        'Q' = 0.0184215914061958
    }
    3 = @{
        'C' = 9652
        'D' = 8568
        'E' = 0.8876916701201824
        'F' = 0.8864059590316573
        'G' = 0.1019541322918984
        'H' = 0.09037275041144063
        'I' = 48077747.45953142
        'J' = 17574124.66950491
        'L' = 17574124.66950491
        'M' = 65651872.12903633
        'N' = 836946297.5090281
        'O' = 819466121.4849579
        'P' = 0.02099791195899919
        'Q' = 0.02144582211362047
    }
    4 = @{
        'C' = 9836
        'D' = 8719
        'E' = 0.8864375762505083
        'F' = 0.8844593223777643
        'G' = 0.1063217547815781
        'H' = 0.09403726718812941
        'I' = 54561128.47940587
        'J' = 20495893.04877832
        'L' = 20495893.04877832
        'M' = 75057021.52818419
        'N' = 875020684.362587
        'O' = 857571736.3566331
        'P' = 0.02342332405971483
        'Q' = 0.02389991668318558
    }
    5 = @{
        'C' = 10032
        'D' = 8902
        'E' = 0.8873604465709729
        'F' = 0.8854187388104238
        'G' = 0.1094883113450708
        'H' = 0.09694300254563555
        'I' = 60667658.19851614
        'J' = 23216084.20014448
        'L' = 23216084.20014448
        'M' = 83883742.3986606
        'N' = 914085108.69052
        'O' = 896599003.2275469
        'P' = 0.02539816476542634
        'Q' = 0.02589349766904937
    }
    6 = @{
        'C' = 10236
        'D' = 9104
        'E' = 0.889409925752247
        'F' = 0.8878486444314414
        'G' = 0.1086977081952903
        'H' = 0.09650711287399286
        'I' = 64007671.86784674
        'J' = 24515046.37498279
        'L' = 24515046.37498279
        'M' = 88522718.24282953
        'N' = 955327879.892617
        'O' = 937736054.0194355
        'P' = 0.02566139530831905
        'Q' = 0.02614280027935738
    }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
